$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 55558892
$ws.Range("I86").Value = 2690.923
$ws.Range("J86").Value = 200005020
$ws.Range("K86").Value = 2690.923
$ws.Range("L86").Value = 200005020
$ws.Range("M86").Value = -1567.923
$ws.Range("N86").Value = -200007266
$ws.Range("H89").Value = 55558892
$ws.Range("I89").Value = 2690.923
$ws.Range("J89").Value = 200005020
$ws.Range("K89").Value = 13454.615
$ws.Range("L89").Value = 1000025100
$ws.Range("M89").Value = -7838.614999999998
$ws.Range("N89").Value = -1000036332
$ws.Range("H132").Value = 1819997.5
$ws.Range("I132").Value = 1833.2554
$ws.Range("J132").Value = 12501713
$ws.Range("K132").Value = 5499.7662
$ws.Range("L132").Value = 37505139
$ws.Range("M132").Value = -2969.7662
$ws.Range("N132").Value = -37510199
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1120.2174
$ws.Range("I45").Value = 1035.3125
$ws.Range("J45").Value = 1314.2858
$ws.Range("K45").Value = 1035.3125
$ws.Range("L45").Value = 1314.2858
$ws.Range("M45").Value = -658.3125
$ws.Range("N45").Value = -2068.2858
$ws.Range("H61").Value = 2546.25
$ws.Range("I61").Value = 2633.5386
$ws.Range("J61").Value = 2168
$ws.Range("K61").Value = 2633.5386
$ws.Range("L61").Value = 2168
$ws.Range("M61").Value = -2421.5386
$ws.Range("N61").Value = -2592
$ws.Range("H97").Value = 2977.5454
$ws.Range("I97").Value = 3883.2666
$ws.Range("J97").Value = 1036.7142
$ws.Range("K97").Value = 3883.2666
$ws.Range("L97").Value = 1036.7142
$ws.Range("M97").Value = -3387.2666
$ws.Range("N97").Value = -2028.7142
$ws.Range("H102").Value = 5233.846
$ws.Range("I102").Value = 3296
$ws.Range("J102").Value = 11693.333
$ws.Range("K102").Value = 3296
$ws.Range("L102").Value = 11693.333
$ws.Range("M102").Value = -1674
$ws.Range("N102").Value = -14937.333
$ws.Range("H122").Value = 1219.7778
$ws.Range("I122").Value = 854
$ws.Range("K122").Value = 2562
$ws.Range("M122").Value = -112
$ws.Range("H136").Value = 2546.25
$ws.Range("I136").Value = 2633.5386
$ws.Range("J136").Value = 2168
$ws.Range("K136").Value = 7900.6158
$ws.Range("L136").Value = 6504
$ws.Range("M136").Value = -5350.6158
$ws.Range("N136").Value = -11604
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 71484136
$ws.Range("I20").Value = 82923.875
$ws.Range("J20").Value = 166685740
$ws.Range("K20").Value = 82923.875
$ws.Range("L20").Value = 166685740
$ws.Range("M20").Value = -82676.875
$ws.Range("N20").Value = -166686234
$ws.Range("H86").Value = 38464536
$ws.Range("I86").Value = 55558144
$ws.Range("J86").Value = 3920
$ws.Range("K86").Value = 55558144
$ws.Range("L86").Value = 3920
$ws.Range("M86").Value = -55557021
$ws.Range("N86").Value = -6166
$ws.Range("H89").Value = 38464536
$ws.Range("I89").Value = 55558144
$ws.Range("J89").Value = 3920
$ws.Range("K89").Value = 277790720
$ws.Range("L89").Value = 19600
$ws.Range("M89").Value = -277785104
$ws.Range("N89").Value = -30832
$ws.Range("H99").Value = 21740758
$ws.Range("I99").Value = 41668384
$ws.Range("J99").Value = 1527.2727
$ws.Range("K99").Value = 41668384
$ws.Range("L99").Value = 1527.2727
$ws.Range("M99").Value = -41666886
$ws.Range("N99").Value = -4523.2727
$ws.Range("H105").Value = 4013.158
$ws.Range("I105").Value = 4576.6665
$ws.Range("K105").Value = 4576.6665
$ws.Range("M105").Value = -2829.6665
$ws.Range("H107").Value = 4858.636
$ws.Range("I107").Value = 860.46875
$ws.Range("J107").Value = 132800
$ws.Range("K107").Value = 860.46875
$ws.Range("L107").Value = 132800
$ws.Range("M107").Value = 1059.53125
$ws.Range("N107").Value = -136640
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 6596
$ws.Range("I23").Value = 8000
$ws.Range("J23").Value = 6245
$ws.Range("K23").Value = 8000
$ws.Range("L23").Value = 6245
$ws.Range("M23").Value = -7760
$ws.Range("N23").Value = -6725
$ws.Range("H27").Value = 6596
$ws.Range("I27").Value = 8000
$ws.Range("J27").Value = 6245
$ws.Range("K27").Value = 8000
$ws.Range("L27").Value = 6245
$ws.Range("M27").Value = -7808
$ws.Range("N27").Value = -6629
$ws.Range("H31").Value = 29085.105
$ws.Range("I31").Value = 2966.3704
$ws.Range("J31").Value = 93194.73
$ws.Range("K31").Value = 2966.3704
$ws.Range("L31").Value = 93194.73
$ws.Range("M31").Value = -2671.3704
$ws.Range("N31").Value = -93784.73
$ws.Range("H34").Value = 29085.105
$ws.Range("I34").Value = 2966.3704
$ws.Range("J34").Value = 93194.73
$ws.Range("K34").Value = 2966.3704
$ws.Range("L34").Value = 93194.73
$ws.Range("M34").Value = -2764.3704
$ws.Range("N34").Value = -93598.73
$ws.Range("H86").Value = 2578.8
$ws.Range("I86").Value = 1988
$ws.Range("K86").Value = 1988
$ws.Range("M86").Value = -865
$ws.Range("H89").Value = 2578.8
$ws.Range("I89").Value = 1988
$ws.Range("K89").Value = 9940
$ws.Range("M89").Value = -4324
$ws.Range("H99").Value = 3867.0667
$ws.Range("I99").Value = 3177.1765
$ws.Range("J99").Value = 4769.231
$ws.Range("K99").Value = 3177.1765
$ws.Range("L99").Value = 4769.231
$ws.Range("M99").Value = -1679.1765
$ws.Range("N99").Value = -7765.231
$ws.Range("H105").Value = 4888.8887
$ws.Range("I105").Value = 5375
$ws.Range("K105").Value = 5375
$ws.Range("M105").Value = -3628
$ws.Range("H126").Value = 3867.0667
$ws.Range("I126").Value = 3177.1765
$ws.Range("J126").Value = 4769.231
$ws.Range("K126").Value = 9531.529500000001
$ws.Range("L126").Value = 14307.693
$ws.Range("M126").Value = -7061.529500000001
$ws.Range("N126").Value = -19247.693
$ws.Range("H132").Value = 1117.909
$ws.Range("I132").Value = 813.1
$ws.Range("J132").Value = 4166
$ws.Range("K132").Value = 2439.3
$ws.Range("L132").Value = 12498
$ws.Range("M132").Value = 90.69999999999982
$ws.Range("N132").Value = -17558
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("H97").Value = 1916.84
$ws.Range("I97").Value = 1882.6666
$ws.Range("J97").Value = 1968.1
$ws.Range("K97").Value = 1882.6666
$ws.Range("L97").Value = 1968.1
$ws.Range("M97").Value = -1386.6666
$ws.Range("N97").Value = -2960.1
$ws.Range("H102").Value = 1413.4286
$ws.Range("I102").Value = 1357.3334
$ws.Range("K102").Value = 1357.3334
$ws.Range("M102").Value = 264.6666
$ws.Range("H122").Value = 1218.6666
$ws.Range("I122").Value = 1208.8889
$ws.Range("J122").Value = 1233.3334
$ws.Range("K122").Value = 3626.6667
$ws.Range("L122").Value = 3700.0002
$ws.Range("M122").Value = -1176.6667
$ws.Range("N122").Value = -8600.0002
$ws.Range("H126").Value = 2520.5454
$ws.Range("I126").Value = 2685.75
$ws.Range("J126").Value = 2080
$ws.Range("K126").Value = 8057.25
$ws.Range("L126").Value = 6240
$ws.Range("M126").Value = -5587.25
$ws.Range("N126").Value = -11180
$ws.Range("H132").Value = 1798.4108
$ws.Range("I132").Value = 1282.1025
$ws.Range("J132").Value = 2982.8823
$ws.Range("K132").Value = 3846.3075
$ws.Range("L132").Value = 8948.6469
$ws.Range("M132").Value = -1316.3075
$ws.Range("N132").Value = -14008.6469
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2642.2
$ws.Range("I7").Value = 2606.7856
$ws.Range("J7").Value = 2724.8333
$ws.Range("K7").Value = 2606.7856
$ws.Range("L7").Value = 2724.8333
$ws.Range("M7").Value = -2494.7856
$ws.Range("N7").Value = -2948.8333
$ws.Range("H40").Value = 2243.6667
$ws.Range("I40").Value = 2171.8572
$ws.Range("J40").Value = 2495
$ws.Range("K40").Value = 2171.8572
$ws.Range("L40").Value = 2495
$ws.Range("M40").Value = -2035.8572
$ws.Range("N40").Value = -2767
$ws.Range("H46").Value = 1749.5
$ws.Range("I46").Value = 1749.5
$ws.Range("K46").Value = 1749.5
$ws.Range("M46").Value = -1561.5
$ws.Range("H93").Value = 2308.1462
$ws.Range("I93").Value = 1951.2142
$ws.Range("J93").Value = 3076.923
$ws.Range("K93").Value = 1951.2142
$ws.Range("L93").Value = 3076.923
$ws.Range("M93").Value = -703.2141999999999
$ws.Range("N93").Value = -5572.923
$ws.Range("H126").Value = 2642.2
$ws.Range("I126").Value = 2606.7856
$ws.Range("J126").Value = 2724.8333
$ws.Range("K126").Value = 7820.3568
$ws.Range("L126").Value = 8174.499899999999
$ws.Range("M126").Value = -5350.3568
$ws.Range("N126").Value = -13114.4999
$ws.Range("H132").Value = 2654.5833
$ws.Range("I132").Value = 1807.9143
$ws.Range("K132").Value = 5423.742899999999
$ws.Range("M132").Value = -2893.742899999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 9000
$ws.Range("J49").Value = 9000
$ws.Range("L49").Value = 9000
$ws.Range("N49").Value = -9460
$ws.Range("H122").Value = 50001420
$ws.Range("I122").Value = 58824956
$ws.Range("J122").Value = 1366.6666
$ws.Range("K122").Value = 176474868
$ws.Range("L122").Value = 4099.9998
$ws.Range("M122").Value = -176472418
$ws.Range("N122").Value = -8999.9998
